$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: "Accounting" label + sample value formatted with a custom accounting
# number format (Indian Rupee accounting style).
$ws.Range("A15").Value = "Accounting"

$ws.Range("B15").Font.Size = 16
$ws.Range("B15").NumberFormat = "_ [$₹-439]\ * #,##0.00_ ;_ [$₹-439]\ * \-#,##0.00_ ;_ [$₹-439]\ * ""-""??_ ;_ @_ "
$ws.Range("B15").Value = 123.4

# Row 16: second accounting example using Excel's built-in Accounting ($) format.
$ws.Range("B16").Font.Size = 16
$ws.Range("B16").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"
$ws.Range("B16").Value = 0

$ws.Range("B16").Select()
